# Auto-generated: apply per-cell value updates from the crypto price refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.480.57"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "1.693.01"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'316.27"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").Value = "'0.3909"
$ws.Range("E7").Value = "  -0.80%  "
$ws.Range("D8").Value = "'0.4054"
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "'1.491"
$ws.Range("E9").Value = "  -2.15%  "
$ws.Range("D10").Value = "'1.001"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("D11").Value = "'52.54"
$ws.Range("E11").Value = "  -2.07%  "
$ws.Range("D12").Value = "'0.08784"
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("D13").Value = "'26.72"
$ws.Range("D14").Value = "'7.508"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").Value = "'8.139"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").Value = "'0.00001348"
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("D17").Value = "1.685.98"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("D18").Value = "'98.05"
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("D19").Value = "'0.07160"
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").Value = "'20.57"
$ws.Range("E20").Value = "  +4.18%  "
$ws.Range("D21").Value = "'7.293"
$ws.Range("E21").Value = "  +2.79%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'14.33"
$ws.Range("E23").Value = "  -0.87%  "
$ws.Range("D24").Value = "24.469.84"
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("D25").Value = "'3.012"
$ws.Range("E25").Value = "  -7.28%  "
$ws.Range("D26").Value = "'2.341"
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("D27").Value = "'22.71"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").Value = "'167.45"
$ws.Range("E28").Value = "  +3.10%  "
$ws.Range("D29").Value = "'8.451"
$ws.Range("E29").Value = "  -4.34%  "
$ws.Range("D30").Value = "'5.386"
$ws.Range("E30").Value = "  +4.00%  "
$ws.Range("D31").Value = "'138.63"
$ws.Range("E31").Value = "  +1.67%  "
$ws.Range("B32").Value = "WEMIXTOKEN"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").Value = "'2.227"
$ws.Range("E32").Value = "  +12.06%  "
$ws.Range("B33").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C33").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D33").Value = "1.872.03"
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("D34").Value = "'0.08748"
$ws.Range("E34").Value = "  -1.94%  "
$ws.Range("D35").Value = "'7.306"
$ws.Range("E35").Value = "  -7.16%  "
$ws.Range("D36").Value = "'1.039"
$ws.Range("E36").Value = "  -4.15%  "
$ws.Range("D37").Value = "'0.02984"
$ws.Range("E37").Value = "  +6.69%  "
$ws.Range("D38").Value = "'0.2789"
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("D39").Value = "'10.95"
$ws.Range("E39").Value = "  -1.46%  "
$ws.Range("D40").Value = "'0.09174"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("D41").Value = "'0.8051"
$ws.Range("E41").Value = "  +4.10%  "
$ws.Range("D42").Value = "'14.18"
$ws.Range("E42").Value = "  -3.02%  "
$ws.Range("D43").Value = "'1.474"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").Value = "'17.58"
$ws.Range("E44").Value = "  +9.71%  "
$ws.Range("D45").Value = "'2.677"
$ws.Range("E45").Value = "  +4.04%  "
$ws.Range("D46").Value = "'0.7270"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("D47").Value = "'4.269"
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("D48").Value = "'1.406"
$ws.Range("E48").Value = "  +4.24%  "
$ws.Range("D49").Value = "'1.002"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").Value = "'140.00"
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").Value = "'0.08173"
$ws.Range("E51").Value = "  +2.37%  "
